$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 475, pushing old rows 475-591 down to 476-592.
$ws.Rows.Item(475).Insert()

# Populate the newly inserted row 475 with the new data record.
$ws.Cells.Item(475, 1).Value2  = 6
$ws.Cells.Item(475, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(475, 3).Value2  = "Metropolitana"
$ws.Cells.Item(475, 4).Value2  = 44782
$ws.Cells.Item(475, 5).Value2  = 13
$ws.Cells.Item(475, 6).Value2  = 100112012
$ws.Cells.Item(475, 7).Value2  = "Espinaca"
$ws.Cells.Item(475, 8).Value2  = "Sin especificar"
$ws.Cells.Item(475, 9).Value2  = "Primera"
$ws.Cells.Item(475, 10).Value2 = 370
$ws.Cells.Item(475, 11).Value2 = 6000
$ws.Cells.Item(475, 12).Value2 = 7000
$ws.Cells.Item(475, 13).Value2 = 6622
$ws.Cells.Item(475, 14).Value2 = "`$/cuna 10 kilos"
$ws.Cells.Item(475, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(475, 16).Value2 = 662
$ws.Cells.Item(475, 17).Value2 = 10
$ws.Cells.Item(475, 18).Value2 = "Hortaliza"
